# Edit script: apply the "demography_eurostat" tutorial update to examples.xlsx
#  - insert a new "immigration" sheet (between "deaths" and "pop_births_deaths")
#  - refresh a handful of population figures for France (2014/2015) that
#    ripple into every sheet that embeds the "pop" table
#  - refresh the narrow-format totals for France (2014/2015)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "immigration" worksheet right after "deaths"
# ---------------------------------------------------------------------------
$deaths = $wb.Worksheets.Item("deaths")
$immigration = $wb.Worksheets.Add($null, $deaths)
$immigration.Name = "immigration"

$immigration.Cells.Item(1,1).Value = "country"
$immigration.Cells.Item(1,2).Value = "citizenship"
$immigration.Cells.Item(1,3).Value = "gender\time"
$immigration.Cells.Item(1,4).Value = 2013
$immigration.Cells.Item(1,5).Value = 2014
$immigration.Cells.Item(1,6).Value = 2015

$data = @(
    @("Belgium",     "Belgium",     "Male",   8822, 10512, 11378),
    @("Belgium",     "Belgium",     "Female", 5727,  6301,  6486),
    @("Belgium",     "Luxembourg",  "Male",    102,   117,   105),
    @("Belgium",     "Luxembourg",  "Female",  117,   123,   114),
    @("Belgium",     "Netherlands", "Male",   4185,  4222,  4183),
    @("Belgium",     "Netherlands", "Female", 3737,  3844,  3942),
    @("Luxembourg",  "Belgium",     "Male",    896,   937,   880),
    @("Luxembourg",  "Belgium",     "Female",  574,   655,   622),
    @("Luxembourg",  "Luxembourg",  "Male",    694,   722,   660),
    @("Luxembourg",  "Luxembourg",  "Female",  607,   586,   535),
    @("Luxembourg",  "Netherlands", "Male",    160,   165,   147),
    @("Luxembourg",  "Netherlands", "Female",   92,    97,    85),
    @("Netherlands", "Belgium",     "Male",   1063,  1141,  1113),
    @("Netherlands", "Belgium",     "Female",  980,  1071,  1181),
    @("Netherlands", "Luxembourg",  "Male",     23,    43,    59),
    @("Netherlands", "Luxembourg",  "Female",   24,    34,    46),
    @("Netherlands", "Netherlands", "Male",  19374, 20037, 21119),
    @("Netherlands", "Netherlands", "Female",16945, 17411, 18084)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $immigration.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}

# ---------------------------------------------------------------------------
# 2) Update the France population figures (2014 / 2015) wherever the "pop"
#    table (country / gender / 2013 / 2014 / 2015) is embedded
# ---------------------------------------------------------------------------
function Update-PopFranceRows($ws, $maleRow, $femaleRow) {
    $ws.Range("D$maleRow").Value = 32045129
    $ws.Range("E$maleRow").Value = 32174258
    $ws.Range("D$femaleRow").Value = 34120851
    $ws.Range("E$femaleRow").Value = 34283895
}

# pop: rows 4 (France/Male) & 5 (France/Female)
Update-PopFranceRows $wb.Worksheets.Item("pop") 4 5

# pop_births_deaths: first block repeats the same "pop" table on rows 4 & 5
Update-PopFranceRows $wb.Worksheets.Item("pop_births_deaths") 4 5

# pop_missing_axis_name: same layout as "pop"
Update-PopFranceRows $wb.Worksheets.Item("pop_missing_axis_name") 4 5

# pop_missing_values: only the France/Female row (row 4) is present
$missingValues = $wb.Worksheets.Item("pop_missing_values")
$missingValues.Range("D4").Value = 34120851
$missingValues.Range("E4").Value = 34283895

# ---------------------------------------------------------------------------
# 3) Update the narrow-format totals for France (2014 / 2015)
# ---------------------------------------------------------------------------
$narrow = $wb.Worksheets.Item("pop_narrow_format")
$narrow.Range("C6").Value = 66165980
$narrow.Range("C7").Value = 66458153

# Keep the originally-active tab selected (it was the last sheet, as before)
$narrow.Activate()
